$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that differ between row 27 and row 28 and need to be swapped.
$cols = @("A", "B", "D", "E", "F", "G", "H", "AC", "AH", "AI")

foreach ($col in $cols) {
    $addr27 = "$col" + "27"
    $addr28 = "$col" + "28"
    $v27 = $ws.Range($addr27).Value2
    $v28 = $ws.Range($addr28).Value2
    $ws.Range($addr27).Value2 = $v28
    $ws.Range($addr28).Value2 = $v27
}
